# DCAF Benchmark.xlsx edit
# - Rename Sheet1 -> 9030_PSP
# - Add a new sheet 9068_IO after it with the 9068 scan-engine-module rate
#   test data (32/64 variable, stock vs modified "Scan Engine" module)
# - Make 9068_IO the active sheet/tab
# - Apply a 0.000000 number format to the one "modified" Duration-Input cell
#   that the author formatted differently (H5)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "9030_PSP"

# Leave sheet1's old "working" selection behind as a range selection (matches
# the header row of its own table) before we switch focus to the new sheet.
$ws1.Range("B1:J1").Select()

# New sheet, inserted immediately after 9030_PSP.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "9068_IO"

# --- Row 1: purpose banner ------------------------------------------------
$ws2.Range("A1").Value = "Purpose"
$ws2.Range("B1").Value = "Compare the generic PSP module with a modification to use the Scanned Access API"

# --- Row 3: table header (mirrors the 9030_PSP sheet's own header) -------
$ws2.Range("A3").Value = "Test:"
$ws2.Range("B3").Value = "Tags"
$ws2.Range("C3").Value = "Rate"
$ws2.Range("D3").Value = "Modules"
$ws2.Range("E3").Value = "CPU"
$ws2.Range("F3").Value = "Timed Structures"
$ws2.Range("G3").Value = "Duration Process"
$ws2.Range("H3").Value = "Duration Input"
$ws2.Range("I3").Value = "Duration Output"
$ws2.Range("J3").Value = "Target"

# --- Rows 4-7 data: A/B/C columns first (tags/count/rate label) ----------
$ws2.Range("A4").Value = "1x9205"
$ws2.Range("B4").Value = 32
$ws2.Range("C4").Value = "1k"

$ws2.Range("A5").Value = "1x9205"
$ws2.Range("B5").Value = 32
$ws2.Range("C5").Value = "1k"

$ws2.Range("A6").Value = "2x9206"
$ws2.Range("B6").Value = 64
$ws2.Range("C6").Value = "1k"

$ws2.Range("A7").Value = "2x9206"
$ws2.Range("B7").Value = 64
$ws2.Range("C7").Value = "1k"

# --- Rows 4-7 data: D/E/F/H/J columns (module + measured values) ---------
$ws2.Range("D4").Value = "Scan Engine"
$ws2.Range("E4").Value = 8.24
$ws2.Range("F4").Value = 3.13
$ws2.Range("H4").Value = 0.00018
$ws2.Range("J4").Value = 9068

$ws2.Range("D5").Value = "Mod - Scan Engine"
$ws2.Range("E5").Value = 7.29
$ws2.Range("F5").Value = 2.22
$ws2.Range("H5").Value = 0.00009
$ws2.Range("J5").Value = 9068

$ws2.Range("D6").Value = "Scan Engine"
$ws2.Range("E6").Value = 14.85
$ws2.Range("F6").Value = 4.89
$ws2.Range("H6").Value = 0.00034
$ws2.Range("J6").Value = 9068

$ws2.Range("D7").Value = "Mod - Scan Engine"
$ws2.Range("E7").Value = 8.22
$ws2.Range("F7").Value = 2.83
$ws2.Range("H7").Value = 0.00014
$ws2.Range("J7").Value = 9068

# --- Row 2: note banner, added last -------------------------------------
$ws2.Range("A2").Value = "Note on values:"
$ws2.Range("B2").Value = "All values are representative snapshots. Not taken through statistical processes"

# The modified-module Duration Input sample (H5) gets a finer fixed-point
# display format than the rest of the column.
$ws2.Range("H5").NumberFormat = "0.000000"

# Column widths tuned for the new sheet's longer labels.
$ws2.Columns.Item(1).ColumnWidth = 16.666666666666668
$ws2.Columns.Item(4).ColumnWidth = 18.166666666666668
$ws2.Columns.Item(6).ColumnWidth = 17.666666666666668
$ws2.Columns.Item(7).ColumnWidth = 19.666666666666668
$ws2.Columns.Item(8).ColumnWidth = 17
$ws2.Columns.Item(9).ColumnWidth = 16.666666666666668

# Leave the same cell selected/active that the author had when they saved.
$ws2.Range("I14").Select()
